# Correcion validacion fecha inventario
#
# Applies to "Base clientes.xlsx" (sheet "Base Clientes carga manual"):
#  - Corrects the Num_Distri for row 5 (DIEGO FRANCO LEAL GARCIA) from 10210128 to 500226
#  - Adds 5 new client rows to the Tabla1 table (rows 11-15)
#  - Leaves the selection on C6, as last left by the editor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Fix an incorrect distributor number on row 5
$ws.Range("A5").Value = 500226

# New clients to append to the table
$newClients = @(
    @{ Num = 62000076;  Name = "DIABONOS S.A." },
    @{ Num = 10236216;  Name = "ARIANNA GARCIA" },
    @{ Num = 10234501;  Name = "Agroquimicos Libra" },
    @{ Num = 10220649;  Name = "ASESORIA INTEGRAL LUMINARIAS" },
    @{ Num = 10234501;  Name = "Agroquimicos Libra" }
)

foreach ($client in $newClients) {
    $row = $tbl.ListRows.Add()
    $rowIndex = $tbl.Range.Rows.Count + $tbl.Range.Row - 1
    $ws.Cells.Item($rowIndex, 1).Value = $client.Num
    $ws.Cells.Item($rowIndex, 2).Value = $client.Name
}

# Leave selection where the editor left it
$ws.Range("C6").Select() | Out-Null
